$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.399.77'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6299'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07677'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2939'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07746'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '1.843.52'
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.00001094'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +9.48%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.012'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6798'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').Value = '2.094.29'
$ws.Range('E17').Value = '  -7.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.140'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '29.419.30'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '229.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.364'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.311'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.41%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.466'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05678'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.048'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('E35').Value = '  +0.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7102'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.780'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = '1.231.12'
$ws.Range('E39').Value = '  -1.92%  '
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.470'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9145'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.38'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.17'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000121'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.71%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.164'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.70%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4013'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.062'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.690'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1124'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.09%  '
